$wb = $excel.ActiveWorkbook

# Add new worksheet '4.7' at the very end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = '4.7'

# Populate the vocabulary table (A:French, B:Dutch, C:status)
$data = @(
  @('une amitié', 'een vriendschap', 'wrong'),
  @('autrefois', 'voorheen, vroeger', 'correct'),
  @('une bière', 'een bier', 'correct'),
  @('blond', 'blond', 'correct'),
  @('calme', 'kalm, rustig', 'correct'),
  @('un camion', 'een vrachtwagen', 'correct'),
  @('une cellule', 'een cel', 'correct'),
  @('chasser', 'jagen [op], verjagen', 'correct'),
  @('une chemise', 'een hemd, een overhemd', 'correct'),
  @('un, une commissaire', 'een commissaris', 'correct'),
  @('se dépêcher', 'zich haasten', 'correct'),
  @('dirigeant', 'leidinggevend, leidend', 'wrong'),
  @('un disque', 'een schijf', 'wrong'),
  @('un dommage', 'een schade', 'correct'),
  @('une émission', 'een programma, een uitzending', 'correct'),
  @('une énergie', 'een energie', 'correct'),
  @('enforcer', 'indrukken, intrappen', 'correct'),
  @('s''enforcer dans', 'wegzinken in, wegzakken in', 'wrong'),
  @('un entraineur', 'een trainer', 'correct'),
  @('un entretien', 'een onderhoud, een interview', 'correct'),
  @('un établissement', 'een instelling', 'wrong'),
  @('examiner', 'onderzoeken', 'correct'),
  @('une fédération', 'een federatie', 'wrong'),
  @('ferme', 'stevig', 'correct'),
  @('un gardien', 'een bewaker', 'correct'),
  @('gêner', 'storen, verstoren, in verlegenheid brengen', 'correct'),
  @('se gêner', 'zich schamen, elkaar in de weg lopen', 'correct'),
  @('une grève', 'een staking', 'wrong'),
  @('incroyable', 'ongeloofelijk, geweldig', 'correct'),
  @('une instruction', 'een opdracht, een instructie', 'correct'),
  @('judiciaire', 'gerechtelijk', 'wrong'),
  @('un lieutenant', 'een luitenant', 'correct'),
  @('une manifestation', 'een betoging', 'wrong'),
  @('une menace', 'een bedreiging', 'wrong'),
  @('la neige', 'de sneeuw', 'correct'),
  @('parfait', 'perfect', 'correct'),
  @('un particulier', 'een particulier, een individu', 'correct'),
  @('un, une pauvre', 'een arme', 'correct'),
  @('une pomme', 'een appel', 'correct'),
  @('un portrait', 'een portret', 'correct'),
  @('une pratique', 'een gebruik', 'wrong'),
  @('un rang', 'een rang', 'correct'),
  @('ranger', 'opruimen', 'correct'),
  @('la retraite', 'het pensioen, het terugtrekken', 'wrong'),
  @('un rythme', 'een ritme', 'correct'),
  @('une séance', 'een zitting', 'wrong'),
  @('le secours', 'de hulp', 'wrong'),
  @('un talent', 'een talent', 'correct'),
  @('téléphoner [à]', 'telefoneren [naar]', 'correct'),
  @('se téléphoner', 'elkaar opbellen', 'correct'),
  @('total', 'totaal', 'correct'),
  @('une urgence', 'een noodsituatie', 'correct'),
  @('vaste', 'uitgebreid, ruim', 'correct')
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 1
  $rowVals = $data[$i]
  for ($j = 0; $j -lt $rowVals.Length; $j++) {
    $ws.Cells.Item($r, $j + 1).Value = $rowVals[$j]
  }
}

# Re-activate sheet '2.5' (index 14, 1-based) as the active tab
$wb.Worksheets.Item(14).Activate()

